# 00 - 04 fig 1 updates
# Appends new figure rows (31-37) to Sheet1, matching the existing table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A original paper, B figure name, C figure number, D year, E page number,
#          F image rotation, G new image path, H subcategory, I is_viewed (boolean TRUE)

$newRows = @(
    @("Bacharach et al._2000_ASQ_Boundary Management Tactics and Logics of Action.pdf", "F1_P15_Bacharach et al._2000_ASQ_Boundary Management Tactics and Logics of Action.png", 1, 2000, 16, 0, "final_figures/2000/F1_P15_Bacharach et al._2000_ASQ_Boundary Management Tactics and Logics of Action.png", "2x2 matrix", $true),
    @("Maznevski & Chudoba_2000_OrgSci_Bridging Qualitative and Quantitative Methods.pdf", "F1_P14_Maznevski & Chudoba_2000_OrgSci_Bridging Qualitative and Quantitative Methods.png", 1, 2000, 15, 0, "final_figures/2000/F1_P14_Maznevski & Chudoba_2000_OrgSci_Bridging Qualitative and Quantitative Methods.png", "process diagram", $true),
    @("Maznevski & Chudoba_2000_OrgSci_Bridging Space over Time.pdf", "F1_P15_Maznevski & Chudoba_2000_OrgSci_Bridging Space over Time.png", 1, 2000, 16, 0, "final_figures/2000/F1_P15_Maznevski & Chudoba_2000_OrgSci_Bridging Space over Time.png", "process diagram", $true),
    @("Orlikowski_2000_OrgSci_Using Technology and Constituting Structures.pdf", "F1_P7_Orlikowski_2000_OrgSci_Using Technology and Constituting Structures.png", 1, 2000, 8, 0, "final_figures/2000/F1_P7_Orlikowski_2000_OrgSci_Using Technology and Constituting Structures.png", "none selected", $true),
    @("Pratt_2000_ASQ_The Good, the Bad and the Ambivalent.pdf", "F1_P21_Pratt_2000_ASQ_The Good, the Bad and the Ambivalent.png", 1, 2000, 22, 0, "final_figures/2000/F1_P21_Pratt_2000_ASQ_The Good, the Bad and the Ambivalent.png", "conceptual diagram", $true),
    @("Shane_2000_OrgSci_Prior Knowledge and the Discovery of Entrepreneurial Opportunities.pdf", "F1_P6_Shane_2000_OrgSci_Prior Knowledge and the Discovery of Entrepreneurial Opportunities.png", 1, 2000, 7, 0, "final_figures/2000/F1_P6_Shane_2000_OrgSci_Prior Knowledge and the Discovery of Entrepreneurial Opportunities.png", "conceptual diagram", $true),
    @("Voss et al._2000_OrgSci_Linking Organizational Values to Relationships with External Constituents_Quant.pdf", "F1_P5_Voss et al._2000_OrgSci_Linking Organizational Values to Relationships with External Constituents_Quant.png", 1, 2000, 6, 0, "final_figures/2000/F1_P5_Voss et al._2000_OrgSci_Linking Organizational Values to Relationships with External Constituents_Quant.png", "conceptual diagram", $true)
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
}
